$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$nl = [char]10

$text = "Conversión del día 💰" + $nl + `
"✅ Dólar paralelo: 68" + $nl + `
$nl + `
"Binance" + $nl + `
"✅ 1000 Bs = 9.07 = 37895.08 pesos" + $nl + `
"✅ 37895.08 pesos = 9.06 = 974.7 Bs" + $nl + `
$nl + `
"Promedio competencia" + $nl + `
"✅ Tasa pesos: 20" + $nl + `
"✅ Tasa Bs: 20" + $nl + `
"✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $text

$wsTasas.Range("N10").Value = 110.199
$wsTasas.Range("O10").Value = 4176
$wsTasas.Range("N12").Value = 4182
$wsTasas.Range("O12").Value = 107.565
